# Append 6 new "particella" rows (rows 32-37) to Sheet1, matching the
# existing table's layout: col A = sequential index (styled like the rest
# of column A), col B = particella code (plain text), col C = comune code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing index cell (A31) so the new
# index cells (A32:A37) pick up the same cell style already used by the
# sheet (bold, bordered, centered) instead of minting a brand-new style.
$ws.Range("A31").Copy()
$ws.Range("A32:A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @(30, "47/3",   277),
    @(31, "302/1",  277),
    @(32, "2129/1", 394),
    @(33, "2129/2", 394),
    @(34, "2129/3", 394),
    @(35, "2103/7", 394)
)

$r = 32
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
